# Insert a new weekly "Haba" price record as row 31 in the historical log,
# shifting all subsequent rows (old 31..102) down by one (new 32..103).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("31:31").Insert()

$ws.Range("A31").Value = 3
$ws.Range("B31").Value = "Femacal de La Calera"
$ws.Range("C31").Value = "Coquimbo"
$ws.Range("D31").Value = 44526
$ws.Range("E31").Value = 5
$ws.Range("F31").Value = 100112026
$ws.Range("G31").Value = "Haba"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 40
$ws.Range("K31").Value = 8000
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = 8000
$ws.Range("N31").Value = "`$/malla 25 kilos"
$ws.Range("O31").Value = "Provincia de Limarí"
$ws.Range("P31").Value = 320
$ws.Range("Q31").Value = 25
$ws.Range("R31").Value = "Hortaliza"
